# Update cryptocurrency price/volume data on Fri Jun 21 06:33:39 UTC 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '64.625.04'
$ws.Cells.Item(2, 5).Value = '  -1.38%  '
$ws.Cells.Item(3, 4).Value = '3.512.00'
$ws.Cells.Item(3, 5).Value = '  -2.04%  '
$ws.Cells.Item(4, 5).Value = '  +0.02%  '
$ws.Cells.Item(5, 4).Value = '585.40'
$ws.Cells.Item(5, 5).Value = '  -2.84%  '
$ws.Cells.Item(6, 4).Value = '132.53'
$ws.Cells.Item(6, 5).Value = '  -2.66%  '
$ws.Cells.Item(7, 4).Value = '3.512.70'
$ws.Cells.Item(7, 5).Value = '  -1.90%  '
$ws.Cells.Item(8, 5).Value = '  +0.02%  '
$ws.Cells.Item(9, 4).Value = '0.489'
$ws.Cells.Item(9, 5).Value = '  -1.32%  '
$ws.Cells.Item(10, 5).Value = '  -0.51%  '
$ws.Cells.Item(11, 4).Value = '7.20'
$ws.Cells.Item(11, 5).Value = '  -0.02%  '
$ws.Cells.Item(12, 4).Value = '0.388'
$ws.Cells.Item(12, 5).Value = '  -0.84%  '
$ws.Cells.Item(13, 4).Value = '4.110.50'
$ws.Cells.Item(13, 5).Value = '  -1.98%  '
$ws.Cells.Item(14, 4).Value = '27.72'
$ws.Cells.Item(14, 5).Value = '  +0.10%  '
$ws.Cells.Item(15, 5).Value = '  -3.11%  '
$ws.Cells.Item(17, 4).Value = '3.512.30'
$ws.Cells.Item(17, 5).Value = '  -2.02%  '
$ws.Cells.Item(18, 4).Value = '64.608.14'
$ws.Cells.Item(18, 5).Value = '  -1.52%  '
$ws.Cells.Item(19, 4).Value = '9.90'
$ws.Cells.Item(20, 4).Value = '14.22'
$ws.Cells.Item(20, 5).Value = '  -2.34%  '
$ws.Cells.Item(21, 4).Value = '5.69'
$ws.Cells.Item(21, 5).Value = '  -3.62%  '
$ws.Cells.Item(22, 4).Value = '390.56'
$ws.Cells.Item(22, 5).Value = '  -1.20%  '
$ws.Cells.Item(23, 4).Value = '0.579'
$ws.Cells.Item(23, 5).Value = '  -1.29%  '
$ws.Cells.Item(24, 4).Value = '3.653.94'
$ws.Cells.Item(24, 5).Value = '  -2.04%  '
$ws.Cells.Item(25, 4).Value = '74.01'
$ws.Cells.Item(25, 5).Value = '  -0.54%  '
$ws.Cells.Item(26, 5).Value = '  +0.02%  '
$ws.Cells.Item(27, 5).Value = '  -4.32%  '
$ws.Cells.Item(28, 4).Value = '1.55'
$ws.Cells.Item(28, 5).Value = '  -6.90%  '
$ws.Cells.Item(29, 2).Value = 'Binance-PegBSC-USD'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Cells.Item(29, 4).Value = '1.00'
$ws.Cells.Item(29, 5).Value = '  -0.02%  '
$ws.Cells.Item(30, 2).Value = 'RenderToken'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(30, 4).Value = '7.44'
$ws.Cells.Item(30, 5).Value = '  -8.23%  '
$ws.Cells.Item(31, 4).Value = '2.25'
$ws.Cells.Item(31, 5).Value = '  -5.32%  '
$ws.Cells.Item(32, 4).Value = '8.21'
$ws.Cells.Item(32, 5).Value = '  -5.61%  '
$ws.Cells.Item(33, 4).Value = '3.518.15'
$ws.Cells.Item(33, 5).Value = '  -1.52%  '
$ws.Cells.Item(34, 5).Value = '  +0.00%  '
$ws.Cells.Item(35, 4).Value = '24.01'
$ws.Cells.Item(35, 5).Value = '  -1.49%  '
$ws.Cells.Item(36, 5).Value = '  -1.33%  '
$ws.Cells.Item(37, 4).Value = '5.28'
$ws.Cells.Item(37, 5).Value = '  -0.76%  '
$ws.Cells.Item(38, 4).Value = '1.58'
$ws.Cells.Item(38, 5).Value = '  -0.75%  '
$ws.Cells.Item(39, 4).Value = '171.20'
$ws.Cells.Item(39, 5).Value = '  +0.06%  '
$ws.Cells.Item(40, 4).Value = '6.95'
$ws.Cells.Item(40, 5).Value = '  -1.57%  '
$ws.Cells.Item(41, 4).Value = '0.0807'
$ws.Cells.Item(41, 5).Value = '  -3.26%  '
$ws.Cells.Item(42, 2).Value = 'EnergySwap'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(42, 4).Value = '26.62'
$ws.Cells.Item(42, 5).Value = '  +2.15%  '
$ws.Cells.Item(43, 2).Value = 'Mantle'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(43, 4).Value = '0.813'
$ws.Cells.Item(43, 5).Value = '  -2.39%  '
$ws.Cells.Item(44, 5).Value = '  +0.05%  '
$ws.Cells.Item(45, 4).Value = '42.07'
$ws.Cells.Item(45, 5).Value = '  -2.90%  '
$ws.Cells.Item(46, 4).Value = '1.21'
$ws.Cells.Item(46, 5).Value = '  -2.76%  '
$ws.Cells.Item(47, 4).Value = '4.40'
$ws.Cells.Item(47, 5).Value = '  -2.56%  '
$ws.Cells.Item(48, 5).Value = '  -3.23%  '
$ws.Cells.Item(49, 4).Value = '2.464.21'
$ws.Cells.Item(49, 5).Value = '  +0.48%  '
$ws.Cells.Item(50, 4).Value = '6.89'
$ws.Cells.Item(50, 5).Value = '  -1.88%  '
$ws.Cells.Item(51, 2).Value = 'SuiNetwork'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Cells.Item(51, 4).Value = '0.893'
$ws.Cells.Item(51, 5).Value = '  +1.76%  '
